$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 74.75
$ws.Range("I12").Value = 50
$ws.Range("K12").Value = 50
$ws.Range("M12").Value = 120
$ws.Range("H62").Value = 175199.83
$ws.Range("I62").Value = 254199.75
$ws.Range("J62").Value = 17200
$ws.Range("K62").Value = 254199.75
$ws.Range("L62").Value = 17200
$ws.Range("M62").Value = -253575.75
$ws.Range("N62").Value = -18448
$ws.Range("H65").Value = 175199.83
$ws.Range("I65").Value = 254199.75
$ws.Range("J65").Value = 17200
$ws.Range("K65").Value = 1270998.75
$ws.Range("L65").Value = 86000
$ws.Range("M65").Value = -1267878.75
$ws.Range("N65").Value = -92240
$ws.Range("H105").Value = 36800
$ws.Range("J105").Value = 36800
$ws.Range("L105").Value = 36800
$ws.Range("N105").Value = -43788
$ws.Range("H107").Value = 784
$ws.Range("I107").Value = 750
$ws.Range("J107").Value = 826.5
$ws.Range("K107").Value = 750
$ws.Range("L107").Value = 826.5
$ws.Range("M107").Value = 1170
$ws.Range("N107").Value = -4666.5
$ws.Range("H132").Value = 2962.5625
$ws.Range("I132").Value = 3100.0715
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 9300.2145
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -6770.2145
$ws.Range("N132").Value = -11060

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 16639.125
$ws.Range("I28").Value = 7552.1665
$ws.Range("J28").Value = 43900
$ws.Range("K28").Value = 7552.1665
$ws.Range("L28").Value = 43900
$ws.Range("M28").Value = -7360.1665
$ws.Range("N28").Value = -44284
$ws.Range("H32").Value = 8149.9375
$ws.Range("I32").Value = 5977.323
$ws.Range("K32").Value = 5977.323
$ws.Range("M32").Value = -5690.323
$ws.Range("H61").Value = 633.64703
$ws.Range("I61").Value = 629.5
$ws.Range("J61").Value = 700
$ws.Range("K61").Value = 629.5
$ws.Range("L61").Value = 700
$ws.Range("M61").Value = -417.5
$ws.Range("N61").Value = -1124
$ws.Range("H74").Value = 1024.625
$ws.Range("I74").Value = 1548.3334
$ws.Range("J74").Value = 500.91666
$ws.Range("K74").Value = 1548.3334
$ws.Range("L74").Value = 500.91666
$ws.Range("M74").Value = -674.3334
$ws.Range("N74").Value = -2248.91666
$ws.Range("H77").Value = 1024.625
$ws.Range("I77").Value = 1548.3334
$ws.Range("J77").Value = 500.91666
$ws.Range("K77").Value = 7741.666999999999
$ws.Range("L77").Value = 2504.5833
$ws.Range("M77").Value = -3373.666999999999
$ws.Range("N77").Value = -11240.5833
$ws.Range("H99").Value = 16639.125
$ws.Range("I99").Value = 7552.1665
$ws.Range("J99").Value = 43900
$ws.Range("K99").Value = 7552.1665
$ws.Range("L99").Value = 43900
$ws.Range("M99").Value = -4557.1665
$ws.Range("N99").Value = -49890
$ws.Range("H136").Value = 633.64703
$ws.Range("I136").Value = 629.5
$ws.Range("J136").Value = 700
$ws.Range("K136").Value = 1888.5
$ws.Range("L136").Value = 2100
$ws.Range("M136").Value = 661.5
$ws.Range("N136").Value = -7200

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 32375.676
$ws.Range("I134").Value = 2278.0833
$ws.Range("J134").Value = 104609.9
$ws.Range("K134").Value = 6834.249899999999
$ws.Range("L134").Value = 313829.7
$ws.Range("M134").Value = -4299.249899999999
$ws.Range("N134").Value = -318899.7

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 217.45454
$ws.Range("I7").Value = 202.22223
$ws.Range("J7").Value = 286
$ws.Range("K7").Value = 202.22223
$ws.Range("L7").Value = 286
$ws.Range("M7").Value = -89.22223
$ws.Range("N7").Value = -512
$ws.Range("H31").Value = 2853.8484
$ws.Range("I31").Value = 3161.8262
$ws.Range("J31").Value = 2145.5
$ws.Range("K31").Value = 3161.8262
$ws.Range("L31").Value = 2145.5
$ws.Range("M31").Value = -2866.8262
$ws.Range("N31").Value = -2735.5
$ws.Range("H34").Value = 2853.8484
$ws.Range("I34").Value = 3161.8262
$ws.Range("J34").Value = 2145.5
$ws.Range("K34").Value = 3161.8262
$ws.Range("L34").Value = 2145.5
$ws.Range("M34").Value = -2959.8262
$ws.Range("N34").Value = -2549.5
$ws.Range("H134").Value = 18520286
$ws.Range("I134").Value = 1787.4166
$ws.Range("J134").Value = 166668270
$ws.Range("K134").Value = 5362.2498
$ws.Range("L134").Value = 500004810
$ws.Range("M134").Value = -2827.2498
$ws.Range("N134").Value = -500009880

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 50.5
$ws.Range("I12").Value = 40.4
$ws.Range("J12").Value = 56.11111
$ws.Range("K12").Value = 121.2
$ws.Range("L12").Value = 168.33333
$ws.Range("M12").Value = 51.80000000000001
$ws.Range("N12").Value = -514.3333299999999

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 9999.5
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 9999.5
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 9999.5
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -10301.5
$ws.Range("H46").Value = 10506.571
$ws.Range("I46").Value = 4925
$ws.Range("J46").Value = 17948.666
$ws.Range("K46").Value = 4925
$ws.Range("L46").Value = 17948.666
$ws.Range("M46").Value = -4769
$ws.Range("N46").Value = -18260.666
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H101").Value = 29885.666
$ws.Range("J101").Value = 29885.666
$ws.Range("L101").Value = 29885.666
$ws.Range("N101").Value = -36375.666
$ws.Range("H104").Value = 34865.332
$ws.Range("J104").Value = 34865.332
$ws.Range("L104").Value = 34865.332
$ws.Range("N104").Value = -41853.332
$ws.Range("H113").Value = 1393.5
$ws.Range("I113").Value = 553.6667
$ws.Range("K113").Value = 553.6667
$ws.Range("M113").Value = 1616.3333
$ws.Range("H132").Value = 3943.7878
$ws.Range("I132").Value = 4100.9546
$ws.Range("J132").Value = 3629.4546
$ws.Range("K132").Value = 12302.8638
$ws.Range("L132").Value = 10888.3638
$ws.Range("M132").Value = -9772.863799999999
$ws.Range("N132").Value = -15948.3638

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 1000
$ws.Range("I18").Value = 1000
$ws.Range("K18").Value = 1000
$ws.Range("M18").Value = -828
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H46").Value = 1996.8334
$ws.Range("I46").Value = 2497.75
$ws.Range("J46").Value = 995
$ws.Range("K46").Value = 2497.75
$ws.Range("L46").Value = 995
$ws.Range("M46").Value = -2309.75
$ws.Range("N46").Value = -1371
$ws.Range("H132").Value = 2396.3215
$ws.Range("I132").Value = 1932.9333
$ws.Range("J132").Value = 2931
$ws.Range("K132").Value = 5798.7999
$ws.Range("L132").Value = 8793
$ws.Range("M132").Value = -3268.7999
$ws.Range("N132").Value = -13853
$ws.Range("H136").Value = 4173.647
$ws.Range("I136").Value = 1387.6666
$ws.Range("K136").Value = 4162.9998
$ws.Range("M136").Value = -1612.9998

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 68949.75
$ws.Range("J46").Value = 68949.75
$ws.Range("L46").Value = 68949.75
$ws.Range("N46").Value = -69411.75
$ws.Range("H54").Value = 10000
$ws.Range("I54").Value = 10000
$ws.Range("K54").Value = 10000
$ws.Range("M54").Value = -9480
$ws.Range("H81").Value = 1544.1177
$ws.Range("I81").Value = 1539.8
$ws.Range("J81").Value = 1545.9166
$ws.Range("K81").Value = 3079.6
$ws.Range("L81").Value = 3091.8332
$ws.Range("M81").Value = -2018.6
$ws.Range("N81").Value = -5213.8332
$ws.Range("H84").Value = 1544.1177
$ws.Range("I84").Value = 1539.8
$ws.Range("J84").Value = 1545.9166
$ws.Range("K84").Value = 15398
$ws.Range("L84").Value = 15459.166
$ws.Range("M84").Value = -10094
$ws.Range("N84").Value = -26067.166
$ws.Range("H132").Value = 1544.8
$ws.Range("I132").Value = 1167.1904
$ws.Range("K132").Value = 3501.5712
$ws.Range("M132").Value = -971.5711999999999
$ws.Range("H134").Value = 68949.75
$ws.Range("J134").Value = 68949.75
$ws.Range("L134").Value = 206849.25
$ws.Range("N134").Value = -211919.25
$ws.Range("H136").Value = 550.2381
$ws.Range("I136").Value = 474.30768
$ws.Range("K136").Value = 1422.92304
$ws.Range("M136").Value = 1127.07696
